$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.292413
$ws.Range("H2").Value = 42.877239
$ws.Range("I2").Value = 0.3214711970063286
$ws.Range("J2").Value = 0.3214711970063286
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.292413
$ws.Range("N2").Value = 42.877239
$ws.Range("O2").Value = 0.3214711970063286
$ws.Range("P2").Value = 0.3214711970063286
$ws.Range("Q2").Value = 204.2730693625689
$ws.Range("R2").Value = 1838.457624263121
$ws.Range("S2").Value = 0.1033437305046818
$ws.Range("T2").Value = 0.1033437305046818

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.292413
$ws.Range("H3").Value = 42.877239
$ws.Range("I3").Value = 0.3214711970063286
$ws.Range("J3").Value = 0.3214711970063286
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.655853
$ws.Range("N3").Value = 7.967559000000001
$ws.Range("O3").Value = 0.05973660591691893
$ws.Range("P3").Value = 0.05973660591691893
$ws.Range("Q3").Value = 37.958547943289
$ws.Range("R3").Value = 341.626931489601
$ws.Range("S3").Value = 0.01920359820920726
$ws.Range("T3").Value = 0.01920359820920726

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.292413
$ws.Range("H4").Value = 42.877239
$ws.Range("I4").Value = 0.3214711970063286
$ws.Range("J4").Value = 0.3214711970063286
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 27.511123
$ws.Range("N4").Value = 82.533369
$ws.Range("O4").Value = 0.6187921970767525
$ws.Range("P4").Value = 0.6187921970767525
$ws.Range("Q4").Value = 393.2003320097989
$ws.Range("R4").Value = 3538.80298808819
$ws.Range("S4").Value = 0.1989238682924396
$ws.Range("T4").Value = 0.1989238682924396

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.655853
$ws.Range("H5").Value = 7.967559000000001
$ws.Range("I5").Value = 0.05973660591691893
$ws.Range("J5").Value = 0.05973660591691893
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.292413
$ws.Range("N5").Value = 42.877239
$ws.Range("O5").Value = 0.3214711970063286
$ws.Range("P5").Value = 0.3214711970063286
$ws.Range("Q5").Value = 37.958547943289
$ws.Range("R5").Value = 341.626931489601
$ws.Range("S5").Value = 0.01920359820920726
$ws.Range("T5").Value = 0.01920359820920726

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.655853
$ws.Range("H6").Value = 7.967559000000001
$ws.Range("I6").Value = 0.05973660591691893
$ws.Range("J6").Value = 0.05973660591691893
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.655853
$ws.Range("N6").Value = 7.967559000000001
$ws.Range("O6").Value = 0.05973660591691893
$ws.Range("P6").Value = 0.05973660591691893
$ws.Range("Q6").Value = 7.053555157609
$ws.Range("R6").Value = 63.48199641848101
$ws.Range("S6").Value = 0.003568462086473274
$ws.Range("T6").Value = 0.003568462086473274

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.655853
$ws.Range("H7").Value = 7.967559000000001
$ws.Range("I7").Value = 0.05973660591691893
$ws.Range("J7").Value = 0.05973660591691893
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 27.511123
$ws.Range("N7").Value = 82.533369
$ws.Range("O7").Value = 0.6187921970767525
$ws.Range("P7").Value = 0.6187921970767525
$ws.Range("Q7").Value = 73.065498552919
$ws.Range("R7").Value = 657.589486976271
$ws.Range("S7").Value = 0.0369645456212384
$ws.Range("T7").Value = 0.0369645456212384

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 27.511123
$ws.Range("H8").Value = 82.533369
$ws.Range("I8").Value = 0.6187921970767525
$ws.Range("J8").Value = 0.6187921970767525
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.292413
$ws.Range("N8").Value = 42.877239
$ws.Range("O8").Value = 0.3214711970063286
$ws.Range("P8").Value = 0.3214711970063286
$ws.Range("Q8").Value = 393.2003320097989
$ws.Range("R8").Value = 3538.80298808819
$ws.Range("S8").Value = 0.1989238682924396
$ws.Range("T8").Value = 0.1989238682924396

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 27.511123
$ws.Range("H9").Value = 82.533369
$ws.Range("I9").Value = 0.6187921970767525
$ws.Range("J9").Value = 0.6187921970767525
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.655853
$ws.Range("N9").Value = 7.967559000000001
$ws.Range("O9").Value = 0.05973660591691893
$ws.Range("P9").Value = 0.05973660591691893
$ws.Range("Q9").Value = 73.065498552919
$ws.Range("R9").Value = 657.589486976271
$ws.Range("S9").Value = 0.0369645456212384
$ws.Range("T9").Value = 0.0369645456212384

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 27.511123
$ws.Range("H10").Value = 82.533369
$ws.Range("I10").Value = 0.6187921970767525
$ws.Range("J10").Value = 0.6187921970767525
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 27.511123
$ws.Range("N10").Value = 82.533369
$ws.Range("O10").Value = 0.6187921970767525
$ws.Range("P10").Value = 0.6187921970767525
$ws.Range("Q10").Value = 756.8618887211288
$ws.Range("R10").Value = 6811.75699849016
$ws.Range("S10").Value = 0.3829037831630745
$ws.Range("T10").Value = 0.3829037831630745
